$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width changed from 15.42578125 to 15.7109375 (OOXML stored width).
# The COM ColumnWidth setter quantizes to a pixel grid (Calibri 11, MDW-based),
# so 15.7109375 itself is not reachable; 14.833333333333332 lands on the nearest
# achievable stored width (15.666666666666666), the closest possible match.
$ws.Columns.Item(1).ColumnWidth = 14.833333333333332

# Update cell values for A1:B32 to the new computed results
$ws.Range("A1").Value = -0.19925705524686776
$ws.Range("B1").Value = 0.19893259219436743
$ws.Range("A2").Value = -0.17682640472112965
$ws.Range("B2").Value = 0.17556712402760155
$ws.Range("A3").Value = -0.096558432142815676
$ws.Range("B3").Value = 0.096379727240471169
$ws.Range("A4").Value = -0.088379727264365826
$ws.Range("B4").Value = 0.088047677239238453
$ws.Range("A5").Value = -0.085047677252513942
$ws.Range("B5").Value = 0.083934116429140282
$ws.Range("A6").Value = -0.029975433368996107
$ws.Range("B6").Value = 0.029704815344976154
$ws.Range("A7").Value = -0.019704815377867835
$ws.Range("B7").Value = 0.019650783628786694
$ws.Range("A8").Value = -0.0096507836622099497
$ws.Range("B8").Value = 0.0095845958999096759
$ws.Range("A9").Value = -0.0075845959147899933
$ws.Range("B9").Value = 0.0075393611009912931
$ws.Range("A10").Value = -0.0055393611161900225
$ws.Range("B10").Value = 0.0055374316793557199
$ws.Range("A11").Value = -0.0025374316969237753
$ws.Range("B11").Value = 0.0025347787318166226
$ws.Range("A12").Value = 0.00096522124940934262
$ws.Range("B12").Value = -0.00097773913459731077
$ws.Range("A13").Value = 0.0044777391158818958
$ws.Range("B13").Value = -0.0044806666931247463
$ws.Range("A14").Value = -0.0090822244949819719
$ws.Range("B14").Value = 0.0090532887038801135
$ws.Range("A15").Value = -0.0080532887167983347
$ws.Range("B15").Value = 0.0080345842961051517
$ws.Range("A16").Value = -0.0060345843116231812
$ws.Range("B16").Value = 0.0060036701847723251
$ws.Range("A17").Value = -0.0040036702005901148
$ws.Range("B17").Value = 0.0039999999793609575
$ws.Range("A18").Value = -0.016106187488041712
$ws.Range("B18").Value = 0.016092054083877372
$ws.Range("A19").Value = -0.012092054093994609
$ws.Range("B19").Value = 0.012017162944699233
$ws.Range("A20").Value = -0.008017162955477275
$ws.Range("B20").Value = 0.0080057317454826915
$ws.Range("A21").Value = -0.0040057317563588768
$ws.Range("B21").Value = 0.0039999999890385496
$ws.Range("A22").Value = -0.075008691897986424
$ws.Range("B22").Value = 0.074545806045755114
$ws.Range("A23").Value = -0.040502859822161952
$ws.Range("B23").Value = 0.040099872718867502
$ws.Range("A24").Value = -0.020099872772454397
$ws.Range("B24").Value = 0.019999999945747859
$ws.Range("A25").Value = -0.051958683073882739
$ws.Range("B25").Value = 0.05192539345898517
$ws.Range("A26").Value = -0.049425393474232138
$ws.Range("B26").Value = 0.04938576360993352
$ws.Range("A27").Value = -0.046885763625664989
$ws.Range("B27").Value = 0.046670117687799273
$ws.Range("A28").Value = -0.044670117704391998
$ws.Range("B28").Value = 0.044535351065057682
$ws.Range("A29").Value = -0.08135466333239183
$ws.Range("B29").Value = 0.081170451566499047
$ws.Range("A30").Value = -0.021170451723431682
$ws.Range("B30").Value = 0.021024125002023553
$ws.Range("A31").Value = -0.014024125034373114
$ws.Range("B31").Value = 0.014001572517891248
$ws.Range("A32").Value = -0.0040015725575184291
$ws.Range("B32").Value = 0.003999999974524826
